# LOB1042.docx restructuring:
# Several paragraphs swap their text content (a reordering of sections in
# the source document), while paragraph structure / styles / run formatting
# stay the same. We edit the text of each paragraph in place, by (1-based)
# paragraph index, so the edits do not interfere with each other.

$d = $word.ActiveDocument
$br = [char]11   # Word's manual line-break char -> serializes as <w:br/>

# --- Paragraph 6: "Objetivos" body (PT) ---------------------------------
$d.Paragraphs.Item(6).Range.Text = "Óptica geométrica e Física. Comprovações experimentais de física moderna."

# --- Paragraph 7: "Objetivos" body (EN, italic) -------------------------
$d.Paragraphs.Item(7).Range.Text = "Geometric and physical optics. Experimental verification of modern physics."

# --- Paragraph 9: "Docente(s) Responsável(eis)" bullet -------------------
$d.Paragraphs.Item(9).Range.Text = "Verificação experimental das Leis da ótica e suas aplicações. Fenômenos físicos relativos à Física Moderna"

# --- Paragraph 11: "Programa resumido" body -> numbered list w/ breaks --
$d.Paragraphs.Item(11).Range.Text = (
    "1) Refração e reflexão." + $br +
    "2) Espelhos planos e esféricos e lentes delgadas." + $br +
    "3) Polarização." + $br +
    "4) Interferência de ondas planas." + $br +
    "5) Difração." + $br +
    "6) Espectroscopia ótica." + $br +
    "7) Determinação da constante de Planck." + $br +
    "8) Radiação de corpo negro."
)

# --- Paragraph 12: italic body (EN) --------------------------------------
$d.Paragraphs.Item(12).Range.Text = "Experimental verification of optical laws and their applications. Physical phenomena related to modern physics."

# --- Paragraph 14: "Programa" numbered list -> single evaluation sentence
$d.Paragraphs.Item(14).Range.Text = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# --- Paragraph 17: "Avaliação" bullet, 3 bold labels + 3 value runs ------
# Locate each bold label via Find, then replace only the (non-bold) value
# range that follows it, up to the next label (or paragraph end). This
# keeps the "Método:"/"Critério:"/"Norma de recuperação:" bold runs intact.
# Each label/boundary is RE-located right before use (via a fresh
# Document.Content Find) because earlier replacements shift character
# offsets for everything downstream - stale Range positions would misfire.

function Find-LabelEnd([string]$label) {
    $rng = $d.Content
    $rng.Find.Execute($label, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $rng.End
}

function Find-LabelStart([string]$label) {
    $rng = $d.Content
    $rng.Find.Execute($label, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $rng.Start
}

$metodoEnd = Find-LabelEnd "Método: "
$criterioStart = Find-LabelStart "Critério: "
$valMetodo = $d.Range($metodoEnd, $criterioStart)
$valMetodo.Text = "NF 5,0." + $br

$criterioEnd = Find-LabelEnd "Critério: "
$normaStart = Find-LabelStart "Norma de recuperação: "
$valCriterio = $d.Range($criterioEnd, $normaStart)
$valCriterio.Text = "(NF+RC)/2 5,0, onde RC é uma prova de recuperação a ser aplicada." + $br

$normaEnd = Find-LabelEnd "Norma de recuperação: "
$endNorma = $d.Paragraphs.Item(17).Range.End
$valNorma = $d.Range($normaEnd, $endNorma)
$valNorma.Text = (
    "Apostilas do Laboratório de Ensino de Física do IFSC/USP." + $br +
    "RESNICK, R.; HALLIDAY, D. Fundamentos de Física. Vol. 4, LTC (2008)." + $br +
    "TIPLER, P.; MOSCA, G. Física para Cientistas e Engenheiros. Vol. 4, LTC (2008)." + $br +
    "SEARS, F. W.; ZEMANSKY, M. W.; YOUNG, H. D.; FREEDMAN, R. A. Física I, Vol. 4, Pearson Addison Wesley (2009)." + $br +
    "JEWETT Jr, John W.; SERWAY, Raymond A. Princípios de Física. Vol. 4, Thomson Pioneira (2008)."
)

# --- Paragraph 19: Bibliography list -> docente name ---------------------
$d.Paragraphs.Item(19).Range.Text = "230696 - Carlos José Todero Peixoto"
